$d = $word.ActiveDocument

# --- Fix the duplicated comma typo: "updated,, the records" -> "updated, the records" ---
$d.Content.Find.Execute(",,", $true, $false, $false, $false, $false, $true, 1, $false, ",", 2) | Out-Null

# --- Insert the two new weekly-report paragraphs after "Week of 4/7 - 4/13" ---
$weekOf47 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Week of 4/7 " + [char]0x2013 + " 4/13") {
        $weekOf47 = $p
        break
    }
}

$weekOf47.Range.InsertParagraphAfter()
$bodyPara = $weekOf47.Next()
$bodyPara.Range.Text = [char]9 + "We started working on the task of saving all of the league data past runtime.  We decided to work in xml and started with writing all of the league data to an xml file.  We spent some time learning xml as a group and then we were able to write to the xml file without too much trouble.  We had some minor trouble with indenting and nesting the different elements, but by the end of the day Tuesday, we were able to write a well-formatted xml file to store league data.  We started briefly on reading the file back and loading the data into our program, but still have much work to do here.  We also had a code review this week.      "

$bodyPara.Range.InsertParagraphAfter()
$weekOf414 = $bodyPara.Next()
$weekOf414.Range.Text = "Week of 4/14 " + [char]0x2013 + " 4/20"
